# Add a "candybars-mini" sheet: a duplicate of "candybars" that only keeps
# the header plus the first 7 candy-bar rows, plus one brand-new row for
# "2 Musketeers" (same stats as "3 Musketeers"). The new sheet becomes the
# active tab.

$wb  = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("candybars")

# Duplicate the whole "candybars" sheet (placed right after it) so that
# column widths / formatting carry over exactly, then rename it.
$src.Copy($null, $src)
$mini = $wb.Worksheets.Item(2)
$mini.Name = "candybars-mini"

# Trim down to header (row 1) + first 7 candy bars (rows 2-8); drop the rest.
$mini.Rows("9:26").Delete()

# Append the new candy bar: "2 Musketeers" (weight 54, chocolate + nougat,
# available in America) - same stats as the existing "3 Musketeers" row.
$mini.Range("A9").Value = "2 Musketeers"
$mini.Range("B9").Value = 54
$mini.Range("C9").Value = 1
$mini.Range("D9").Value = 0
$mini.Range("E9").Value = 0
$mini.Range("F9").Value = 1
$mini.Range("G9").Value = 0
$mini.Range("H9").Value = 0
$mini.Range("I9").Value = 0
$mini.Range("J9").Value = 0
$mini.Range("K9").Value = "America"

# Page setup / selection on the new sheet.
$mini.PageSetup.Orientation = 1
[void]$mini.Range("C12").Select()

# Make the new sheet the active tab.
$mini.Activate()
